# Auto-generated edit script applying numeric corrections to Garuda_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, LTW, WVR) as captured in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1299.8823
$ws.Range("I28").Value = 205.55556
$ws.Range("J28").Value = 2531
$ws.Range("K28").Value = 205.55556
$ws.Range("L28").Value = 2531
$ws.Range("M28").Value = 279.44444
$ws.Range("N28").Value = -3501
# Row 62
$ws.Range("H62").Value = 2924.5557
$ws.Range("I62").Value = 3040.625
$ws.Range("J62").Value = 1996
$ws.Range("K62").Value = 3040.625
$ws.Range("L62").Value = 1996
$ws.Range("M62").Value = -2416.625
$ws.Range("N62").Value = -3244
# Row 65
$ws.Range("H65").Value = 2924.5557
$ws.Range("I65").Value = 3040.625
$ws.Range("J65").Value = 1996
$ws.Range("K65").Value = 15203.125
$ws.Range("L65").Value = 9980
$ws.Range("M65").Value = -12083.125
$ws.Range("N65").Value = -16220
# Row 96
$ws.Range("H96").Value = 3402.3333
$ws.Range("I96").Value = 3540.3333
$ws.Range("J96").Value = 3333.3333
$ws.Range("K96").Value = 10620.9999
$ws.Range("L96").Value = 9999.999899999999
$ws.Range("M96").Value = -9247.999899999999
$ws.Range("N96").Value = -12745.9999
# Row 132
$ws.Range("H132").Value = 8405655
$ws.Range("I132").Value = 8405655
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 25216965
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -25214435
$ws.Range("N132").ClearContents()
# Row 137
$ws.Range("H137").Value = 1491.5834
$ws.Range("I137").Value = 1411.7778
$ws.Range("J137").Value = 1731
$ws.Range("K137").Value = 4235.3334
$ws.Range("L137").Value = 5193
$ws.Range("M137").Value = -1685.3334
$ws.Range("N137").Value = -10293

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2154.7036
$ws.Range("I61").Value = 1067.3125
$ws.Range("J61").Value = 3736.3635
$ws.Range("K61").Value = 1067.3125
$ws.Range("L61").Value = 3736.3635
$ws.Range("M61").Value = -855.3125
$ws.Range("N61").Value = -4160.363499999999
# Row 74
$ws.Range("H74").Value = 1457.3334
$ws.Range("I74").Value = 1471.0625
$ws.Range("J74").Value = 1347.5
$ws.Range("K74").Value = 1471.0625
$ws.Range("L74").Value = 1347.5
$ws.Range("M74").Value = -597.0625
$ws.Range("N74").Value = -3095.5
# Row 77
$ws.Range("H77").Value = 1457.3334
$ws.Range("I77").Value = 1471.0625
$ws.Range("J77").Value = 1347.5
$ws.Range("K77").Value = 7355.3125
$ws.Range("L77").Value = 6737.5
$ws.Range("M77").Value = -2987.3125
$ws.Range("N77").Value = -15473.5
# Row 136
$ws.Range("H136").Value = 2154.7036
$ws.Range("I136").Value = 1067.3125
$ws.Range("J136").Value = 3736.3635
$ws.Range("K136").Value = 3201.9375
$ws.Range("L136").Value = 11209.0905
$ws.Range("M136").Value = -651.9375
$ws.Range("N136").Value = -16309.0905

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1638.9412
$ws.Range("I86").Value = 1336.5714
$ws.Range("K86").Value = 1336.5714
$ws.Range("M86").Value = -213.5714
# Row 89
$ws.Range("H89").Value = 1638.9412
$ws.Range("I89").Value = 1336.5714
$ws.Range("K89").Value = 6682.857
$ws.Range("M89").Value = -1066.857
# Row 107
$ws.Range("H107").Value = 2335.1667
$ws.Range("I107").Value = 2402.2
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 2402.2
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -482.1999999999998
$ws.Range("N107").Value = -5840
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1006.39026
$ws.Range("I58").Value = 670.53845
$ws.Range("J58").Value = 1588.5333
$ws.Range("K58").Value = 670.53845
$ws.Range("L58").Value = 1588.5333
$ws.Range("M58").Value = -467.53845
$ws.Range("N58").Value = -1994.5333
# Row 134
$ws.Range("H134").Value = 1092.5161
$ws.Range("I134").Value = 1075.2858
$ws.Range("J134").Value = 1128.7
$ws.Range("K134").Value = 3225.8574
$ws.Range("L134").Value = 3386.1
$ws.Range("M134").Value = -690.8574000000003
$ws.Range("N134").Value = -8456.1
# Row 136
$ws.Range("H136").Value = 1006.39026
$ws.Range("I136").Value = 670.53845
$ws.Range("J136").Value = 1588.5333
$ws.Range("K136").Value = 2011.61535
$ws.Range("L136").Value = 4765.5999
$ws.Range("M136").Value = 538.38465
$ws.Range("N136").Value = -9865.599900000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1093.1428
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1093.1428
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3279.4284
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -3503.4284
# Row 68
$ws.Range("H68").Value = 471.73685
$ws.Range("I68").Value = 443.9
$ws.Range("J68").Value = 502.66666
$ws.Range("K68").Value = 1331.7
$ws.Range("L68").Value = 1507.99998
$ws.Range("M68").Value = -520.6999999999998
$ws.Range("N68").Value = -3129.99998
# Row 71
$ws.Range("H71").Value = 471.73685
$ws.Range("I71").Value = 443.9
$ws.Range("J71").Value = 502.66666
$ws.Range("K71").Value = 3995.1
$ws.Range("L71").Value = 4523.99994
$ws.Range("M71").Value = 60.90000000000009
$ws.Range("N71").Value = -12635.99994
# Row 109
$ws.Range("H109").Value = 983.3333
$ws.Range("I109").Value = 900
$ws.Range("J109").Value = 1000
$ws.Range("K109").Value = 2700
$ws.Range("L109").Value = 3000
$ws.Range("M109").Value = -1660
$ws.Range("N109").Value = -5080
# Row 131
$ws.Range("H131").Value = 5348.3335
$ws.Range("J131").Value = 960
$ws.Range("L131").Value = 2880
$ws.Range("N131").Value = -12960
# Row 135
$ws.Range("H135").Value = 1093.1428
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1093.1428
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 9838.2852
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -14908.2852
# Row 136
$ws.Range("H136").Value = 2020.1852
$ws.Range("I136").Value = 1696.6666
$ws.Range("J136").Value = 2039.2157
$ws.Range("K136").Value = 5089.9998
$ws.Range("L136").Value = 6117.6471
$ws.Range("M136").Value = 10.0002000000004
$ws.Range("N136").Value = -16317.6471

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1507.4546
$ws.Range("I22").Value = 1635.625
$ws.Range("J22").Value = 1165.6666
$ws.Range("K22").Value = 1635.625
$ws.Range("L22").Value = 1165.6666
$ws.Range("M22").Value = -1340.625
$ws.Range("N22").Value = -1755.6666
# Row 27
$ws.Range("H27").Value = 1507.4546
$ws.Range("I27").Value = 1635.625
$ws.Range("J27").Value = 1165.6666
$ws.Range("K27").Value = 1635.625
$ws.Range("L27").Value = 1165.6666
$ws.Range("M27").Value = -1528.625
$ws.Range("N27").Value = -1379.6666
# Row 132
$ws.Range("H132").Value = 3683.5
$ws.Range("I132").Value = 4021.6924
$ws.Range("J132").Value = 3055.4285
$ws.Range("K132").Value = 12065.0772
$ws.Range("L132").Value = 9166.2855
$ws.Range("M132").Value = -9535.0772
$ws.Range("N132").Value = -14226.2855

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 6007.846
$ws.Range("I136").Value = 6565.174
$ws.Range("K136").Value = 19695.522
$ws.Range("M136").Value = -17145.522
